{"js": "// Apply the \"ci-sandbox-pipeline\" -> \"simple-ci-pipeline\" rewrite.\n// Strategy: work from the bottom of the document upward so that earlier\n// paragraph indices stay valid while later ones are mutated/removed.\n\nconst body = context.document.body;\nlet paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\n// ---- Bottom block (second YAML doc: PipelineRun) ----\n\n// \"    name: ci-sandbox-pipeline\" (pipelineRef.name) -> \"simple-ci-pipeline\"\nparas.items[62].insertText(\"    name: simple-ci-pipeline\", Word.InsertLocation.replace);\n\n// Insert the new \"timeouts:\" section right after that paragraph.\nparas.items[62].insertParagraph(\"  timeouts:\", Word.InsertLocation.after);\nawait context.sync();\n\n// Re-fetch so we can anchor the second inserted line after the first.\nparas = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\nparas.items[63].insertParagraph(\"    pipeline: 15m\", Word.InsertLocation.after);\nawait context.sync();\n\n// \"apiVersion: tekton.dev/v1beta1\" (second doc header) -> \".../v1\"\n// Do a surgical replace of just the trailing \"/v1beta1\" -> \"/v1\" run so we\n// don't disturb the proofErr wrapping around \"apiVersion\"/\"tekton.dev\".\n// There are two \"/v1beta1\" occurrences (first doc + second doc); only the\n// second (PipelineRun) one changes.\n{\n  const hits = body.search(\"/v1beta1\", { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n  hits.items[hits.items.length - 1].insertText(\"/v1\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Remove 14 of the trailing blank paragraphs (keep the first 12 of the 26).\nparas = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\nfor (let i = 55; i >= 42; i--) {\n  paras.items[i].delete();\n}\nawait context.sync();\n\n// Remove the \"#!/usr/bin/env bash\" line.\nparas = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\nparas.items[28].delete();\nawait context.sync();\n\n// Remove the old standalone \"install-requirements\" task\n// (name, runAfter, taskSpec:, steps:) now folded into the first task.\nparas = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\nfor (let i = 23; i >= 20; i--) {\n  paras.items[i].delete();\n}\nawait context.sync();\n\n// Remove the blank paragraph that used to separate the two tasks.\nparas = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\nparas.items[19].delete();\nawait context.sync();\n\n// Remove the git-clone task's workspaces/params block (now unused).\nparas = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\nfor (let i = 18; i >= 11; i--) {\n  paras.items[i].delete();\n}\nawait context.sync();\n\n// \"        name: git-clone\" -> \"        steps:\"\nparas = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\nparas.items[10].insertText(\"        steps:\", Word.InsertLocation.replace);\n\n// \"taskRef:\" -> \"taskSpec:\" (surgical word replace keeps the existing\n// spellStart/spellEnd proofErr wrapper intact, same as the rest of the\n// paragraph's runs).\n{\n  const hits = body.search(\"taskRef\", { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n  hits.items[0].insertText(\"taskSpec\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// \"    - name: fetch-repo\" -> \"    - name: install-requirements\"\nparas.items[8].insertText(\"    - name: install-requirements\", Word.InsertLocation.replace);\n\n// \"  name: ci-sandbox-pipeline\" (first doc's metadata.name) -> \"simple-ci-pipeline\"\nparas.items[3].insertText(\"  name: simple-ci-pipeline\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Apply the \"ci-sandbox-pipeline\" -> \"simple-ci-pipeline\" rewrite.\n# Strategy: work from the bottom of the document upward so that earlier\n# paragraph indices stay valid while later ones are mutated/removed.\n# (Word COM paragraph indices are 1-based.)\n\n$d = $word.ActiveDocument\n\n# ---- Bottom block (second YAML doc: PipelineRun) ----\n\n# \"    name: ci-sandbox-pipeline\" (pipelineRef.name) -> \"simple-ci-pipeline\"\n$p = $d.Paragraphs.Item(63)\n$p.Range.Text = \"    name: simple-ci-pipeline\"\n\n# Insert the new \"timeouts:\" section right after that paragraph.\n$p = $d.Paragraphs.Item(63)\n$p.Range.InsertParagraphAfter()\n$d.Paragraphs.Item(64).Range.Text = \"  timeouts:\"\n$d.Paragraphs.Item(64).Range.InsertParagraphAfter()\n$d.Paragraphs.Item(65).Range.Text = \"    pipeline: 15m\"\n\n# \"apiVersion: tekton.dev/v1beta1\" (second doc header) -> \".../v1\"\n# Scope the Find to just that paragraph's range so the first YAML doc's\n# identical \"/v1beta1\" text (paragraph 1) is left untouched, and so the\n# proofErr wrapping around \"apiVersion\"/\"tekton.dev\" in that paragraph is\n# preserved (only the trailing run's text actually changes).\n$p = $d.Paragraphs.Item(57)\n$r = $p.Range\n$null = $r.Find.Execute(\"/v1beta1\", $false, $false, $false, $false, $false, $true, 1, $false, \"/v1\", 2)\n\n# Remove 14 of the trailing blank paragraphs (keep the first 12 of the 26).\nfor ($i = 56; $i -ge 43; $i--) {\n    $d.Paragraphs.Item($i).Range.Delete()\n}\n\n# Remove the \"#!/usr/bin/env bash\" line.\n$d.Paragraphs.Item(29).Range.Delete()\n\n# Remove the old standalone \"install-requirements\" task\n# (name, runAfter, taskSpec:, steps:) now folded into the first task.\n$startPara = $d.Paragraphs.Item(21)\n$endPara = $d.Paragraphs.Item(24)\n$rng = $d.Range($startPara.Range.Start, $endPara.Range.End)\n$rng.Delete()\n\n# Remove the blank paragraph that used to separate the two tasks.\n$d.Paragraphs.Item(20).Range.Delete()\n\n# Remove the git-clone task's workspaces/params block (now unused).\n$startPara = $d.Paragraphs.Item(12)\n$endPara = $d.Paragraphs.Item(19)\n$rng = $d.Range($startPara.Range.Start, $endPara.Range.End)\n$rng.Delete()\n\n# \"        name: git-clone\" -> \"        steps:\"\n$d.Paragraphs.Item(11).Range.Text = \"        steps:\"\n\n# \"taskRef:\" -> \"taskSpec:\" (scoped find keeps the existing\n# spellStart/spellEnd proofErr wrapper intact, same as the rest of the\n# paragraph's runs).\n$p = $d.Paragraphs.Item(10)\n$r = $p.Range\n$null = $r.Find.Execute(\"taskRef\", $false, $false, $false, $false, $false, $true, 1, $false, \"taskSpec\", 2)\n\n# \"    - name: fetch-repo\" -> \"    - name: install-requirements\"\n$d.Paragraphs.Item(9).Range.Text = \"    - name: install-requirements\"\n\n# \"  name: ci-sandbox-pipeline\" (first doc's metadata.name) -> \"simple-ci-pipeline\"\n$d.Paragraphs.Item(4).Range.Text = \"  name: simple-ci-pipeline\"\n"}
